# Insert a new data row before existing row 121, shifting rows 121-218 down to
# 122-219. Populate the new row 121 with the new weekly data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121 (existing row 121 and below shift down by one)
$ws.Rows.Item(121).Insert()

# Populate the new row 121
$ws.Cells.Item(121, 1).Value = 7
$ws.Cells.Item(121, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(121, 3).Value = "Ñuble"
$ws.Cells.Item(121, 4).Value = 44574
$ws.Cells.Item(121, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(121, 5).Value = 16
$ws.Cells.Item(121, 6).Value = 100112023
$ws.Cells.Item(121, 7).Value = "Brócoli"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 400
$ws.Cells.Item(121, 11).Value = 600
$ws.Cells.Item(121, 12).Value = 650
$ws.Cells.Item(121, 13).Value = 625
$ws.Cells.Item(121, 14).Value = "$/unidad"
$ws.Cells.Item(121, 15).Value = "Región del Maule"
$ws.Cells.Item(121, 16).Value = 625
$ws.Cells.Item(121, 17).Value = 1
$ws.Cells.Item(121, 18).Value = "Hortaliza"
